$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.639.54"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.75%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.505.22"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.21%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.47"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.95"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +5.36%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.502.45"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.590"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +5.79%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.65%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.47%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.110.64"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.21%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.37"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.34%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.669.66"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.496.39"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.46%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.09"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "391.16"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.00"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.98%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.56%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +5.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.25"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +7.22%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.37"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +5.04%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +6.76%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.61%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.64%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +5.03%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +6.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "162.93"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.885"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.73%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.94%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.65%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.60"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.93%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "27.03"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.33%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.818.02"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "43.06"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.42%  "
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "354.85"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.82%  "
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.52"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.43%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.69"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +12.27%  "

Write-Output "Applied 83 cell updates"
